$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Copy the header style from A1 (existing header) onto the new J1 header cell,
# then set its value - this gives J1 the same style index (s="1") as the
# other header cells.
$ws2.Range("A1").Copy()
$ws2.Range("J1").PasteSpecial(-4122)
$ws2.Range("J1").Value = "postal_Address"

# New data cell - leading apostrophe forces a text/quotePrefix cell (matches
# the "Check"-style text entries elsewhere in the sheet).
$ws2.Range("J2").Value = "'2-6 MAWSON"

# Widen the new column to fit the address text.
$ws2.Columns.Item(10).ColumnWidth = 45

# Make "signIn" (sheet2) the active sheet/tab and select the newly added cell.
$ws2.Activate()
$ws2.Range("J2").Select()
